$wb = $excel.ActiveWorkbook

# Strip the "-SOLUTION" suffix from the view-name labels on the
# "View_Name" sheet (A1:A12) -- these were marked as solved answers and
# are being reverted to the plain, un-suffixed names.
$viewSheet = $wb.Worksheets.Item("View_Name")
for ($r = 1; $r -le 12; $r++) {
    $cell = $viewSheet.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 -replace "-SOLUTION$", ""
}

# Move the active tab / selection from "Z-Min" back to "View_Name",
# keeping Z-Min's own last selection (C13) intact.
$zMinSheet = $wb.Worksheets.Item("Z-Min")
$null = $zMinSheet.Range("C13").Select()

$null = $viewSheet.Activate()
$null = $viewSheet.Range("A10").Select()
